$d = $word.ActiveDocument

# --- Change 1: fix punctuation in the "多云" (2nd) paragraph ---
$d.Content.Find.Execute(
    "多云，今天是六一儿童节，又是开心的一天呢.2022年6月2日星期四。",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "多云，今天是六一儿童节，又是开心的一天呢。2022年6月2日星期四,", 2) | Out-Null

# --- Change 2: rewrite the last paragraph (端午节) and append a new
#     paragraph about Gaokao, keeping the _GoBack bookmark at the very
#     end of the (new) last paragraph. We replace the whole last
#     paragraph (including its end-of-paragraph mark) with raw WordML
#     so that paragraph-mark formatting (rFonts hint=eastAsia) and the
#     bookmark placement come out exactly right. ---

$targetPara = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*中国传统端午节*") {
        $targetPara = $p
    }
}
if ($targetPara -eq $null) {
    $targetPara = $d.Paragraphs.Last
}
$target = $d.Range($targetPara.Range.Start, $targetPara.Range.End)

$newXml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>中雨，今天是农历五月初五，是中国传统节日:端午节，这一天我们要吃粽子，赛龙舟。2022年6月7日星期二,</w:t></w:r>
</w:p>
<w:p>
<w:pPr><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr></w:pPr>
<w:r><w:rPr><w:rFonts w:hint="eastAsia"/></w:rPr><w:t>晴，今天是高考第一天，上午考语文，下午考数学。</w:t></w:r>
<w:bookmarkStart w:id="0" w:name="_GoBack"/>
<w:bookmarkEnd w:id="0"/>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
'@

$target.InsertXML($newXml)

Write-Output "ok"
